$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830
$ws.Range("H15").Value = 1732.4
$ws.Range("I15").Value = 1732.4
$ws.Range("K15").Value = 5197.200000000001
$ws.Range("M15").Value = -5028.200000000001
$ws.Range("H55").Value = 316.5
$ws.Range("J55").Value = 500
$ws.Range("L55").Value = 500
$ws.Range("N55").Value = -928
$ws.Range("H58").Value = 1150.7858
$ws.Range("I58").Value = 609.4
$ws.Range("K58").Value = 1828.2
$ws.Range("M58").Value = -1678.2
$ws.Range("H61").Value = 422.5
$ws.Range("I61").Value = 396.66666
$ws.Range("K61").Value = 1189.99998
$ws.Range("M61").Value = -1017.99998
$ws.Range("H76").Value = 3105.2632
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 5000
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 3105.2632
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 5000
$ws.Range("N79").Value = -7184
$ws.Range("H82").Value = 2133
$ws.Range("I82").Value = 2133
$ws.Range("K82").Value = 6399
$ws.Range("M82").Value = -5993
$ws.Range("H85").Value = 2133
$ws.Range("I85").Value = 2133
$ws.Range("K85").Value = 6399
$ws.Range("M85").Value = -4995
$ws.Range("H94").Value = 995
$ws.Range("I94").Value = 995
$ws.Range("K94").Value = 995
$ws.Range("M94").Value = -544
$ws.Range("H96").Value = 625817.4
$ws.Range("I96").Value = 909546.4
$ws.Range("J96").Value = 1613.6
$ws.Range("K96").Value = 2728639.2
$ws.Range("L96").Value = 4840.799999999999
$ws.Range("M96").Value = -2727266.2
$ws.Range("N96").Value = -7586.799999999999
$ws.Range("H101").Value = 352.57144
$ws.Range("I101").Value = 383.16666
$ws.Range("J101").Value = 169
$ws.Range("K101").Value = 1149.49998
$ws.Range("L101").Value = 507
$ws.Range("M101").Value = 472.5000199999999
$ws.Range("N101").Value = -3751
$ws.Range("H104").Value = 1705.625
$ws.Range("I104").Value = 1935.1428
$ws.Range("J104").Value = 99
$ws.Range("K104").Value = 5805.428400000001
$ws.Range("L104").Value = 297
$ws.Range("M104").Value = -4058.428400000001
$ws.Range("N104").Value = -3791
$ws.Range("H107").Value = 854.93335
$ws.Range("I107").Value = 837.7273
$ws.Range("K107").Value = 837.7273
$ws.Range("M107").Value = 1082.2727
$ws.Range("H132").Value = 2702.7273
$ws.Range("I132").Value = 2122.625
$ws.Range("K132").Value = 6367.875
$ws.Range("M132").Value = -3837.875
$ws.Range("H137").Value = 1932.65
$ws.Range("I137").Value = 1816.5
$ws.Range("J137").Value = 2397.25
$ws.Range("K137").Value = 5449.5
$ws.Range("L137").Value = 7191.75
$ws.Range("M137").Value = -2899.5
$ws.Range("N137").Value = -12291.75
$ws.Range("H138").Value = 1922.0851
$ws.Range("I138").Value = 1090.2972
$ws.Range("J138").Value = 4999.7
$ws.Range("K138").Value = 3270.8916
$ws.Range("L138").Value = 14999.1
$ws.Range("M138").Value = 1869.1084
$ws.Range("N138").Value = -25279.1

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1570.4783
$ws.Range("I2").Value = 531.75
$ws.Range("J2").Value = 3944.7144
$ws.Range("K2").Value = 531.75
$ws.Range("L2").Value = 3944.7144
$ws.Range("M2").Value = -418.75
$ws.Range("N2").Value = -4170.7144
$ws.Range("H5").Value = 142.44444
$ws.Range("I5").Value = 55.333332
$ws.Range("K5").Value = 55.333332
$ws.Range("M5").Value = 56.666668
$ws.Range("H61").Value = 4230.4736
$ws.Range("I61").Value = 4230.4736
$ws.Range("K61").Value = 4230.4736
$ws.Range("M61").Value = -4018.4736
$ws.Range("H97").Value = 948.8889
$ws.Range("I97").Value = 927.4
$ws.Range("J97").Value = 975.75
$ws.Range("K97").Value = 927.4
$ws.Range("L97").Value = 975.75
$ws.Range("M97").Value = -431.4
$ws.Range("N97").Value = -1967.75
$ws.Range("H102").Value = 3823.75
$ws.Range("I102").Value = 2465.7727
$ws.Range("K102").Value = 2465.7727
$ws.Range("M102").Value = -843.7727
$ws.Range("H116").Value = 1570.4783
$ws.Range("I116").Value = 531.75
$ws.Range("J116").Value = 3944.7144
$ws.Range("K116").Value = 531.75
$ws.Range("L116").Value = 3944.7144
$ws.Range("M116").Value = 1762.25
$ws.Range("N116").Value = -8532.714400000001
$ws.Range("H122").Value = 7938342
$ws.Range("I122").Value = 9261056
$ws.Range("J122").Value = 2055.5
$ws.Range("K122").Value = 27783168
$ws.Range("L122").Value = 6166.5
$ws.Range("M122").Value = -27780718
$ws.Range("N122").Value = -11066.5
$ws.Range("H132").Value = 1898.2174
$ws.Range("I132").Value = 1877.8182
$ws.Range("K132").Value = 5633.4546
$ws.Range("M132").Value = -3103.4546
$ws.Range("H136").Value = 4230.4736
$ws.Range("I136").Value = 4230.4736
$ws.Range("K136").Value = 12691.4208
$ws.Range("M136").Value = -10141.4208

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1570.4783
$ws.Range("I3").Value = 531.75
$ws.Range("J3").Value = 3944.7144
$ws.Range("K3").Value = 531.75
$ws.Range("L3").Value = 3944.7144
$ws.Range("M3").Value = -417.75
$ws.Range("N3").Value = -4172.7144
$ws.Range("H4").Value = 142.44444
$ws.Range("I4").Value = 55.333332
$ws.Range("K4").Value = 55.333332
$ws.Range("M4").Value = 59.666668
$ws.Range("H30").Value = 1500
$ws.Range("I30").Value = 1500
$ws.Range("K30").Value = 1500
$ws.Range("M30").Value = -1375
$ws.Range("H68").Value = 62500.5
$ws.Range("J68").Value = 62500.5
$ws.Range("L68").Value = 62500.5
$ws.Range("N68").Value = -64122.5
$ws.Range("H69").Value = 75000
$ws.Range("J69").Value = 75000
$ws.Range("L69").Value = 75000
$ws.Range("N69").Value = -76622
$ws.Range("H71").Value = 62500.5
$ws.Range("J71").Value = 62500.5
$ws.Range("L71").Value = 187501.5
$ws.Range("N71").Value = -195613.5
$ws.Range("H72").Value = 75000
$ws.Range("J72").Value = 75000
$ws.Range("L72").Value = 225000
$ws.Range("N72").Value = -233112
$ws.Range("H86").Value = 22567520
$ws.Range("I86").Value = 2908.5715
$ws.Range("J86").Value = 259495950
$ws.Range("K86").Value = 2908.5715
$ws.Range("L86").Value = 259495950
$ws.Range("M86").Value = -1785.5715
$ws.Range("N86").Value = -259498196
$ws.Range("H89").Value = 22567520
$ws.Range("I89").Value = 2908.5715
$ws.Range("J89").Value = 259495950
$ws.Range("K89").Value = 14542.8575
$ws.Range("L89").Value = 1297479750
$ws.Range("M89").Value = -8926.8575
$ws.Range("N89").Value = -1297490982
$ws.Range("H105").Value = 2580.5715
$ws.Range("I105").Value = 2204.6538
$ws.Range("K105").Value = 2204.6538
$ws.Range("M105").Value = -457.6538

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2921.1
$ws.Range("I31").Value = 1230.8823
$ws.Range("K31").Value = 1230.8823
$ws.Range("M31").Value = -935.8823
$ws.Range("H34").Value = 2921.1
$ws.Range("I34").Value = 1230.8823
$ws.Range("K34").Value = 1230.8823
$ws.Range("M34").Value = -1028.8823
$ws.Range("H55").Value = 29999.25
$ws.Range("J55").Value = 29999.25
$ws.Range("L55").Value = 29999.25
$ws.Range("N55").Value = -30629.25
$ws.Range("H58").Value = 2266.2173
$ws.Range("I58").Value = 2266.2173
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2266.2173
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -2063.2173
$ws.Range("H62").Value = 6968.3335
$ws.Range("I62").Value = 6968.3335
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 6968.3335
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6344.3335
$ws.Range("H65").Value = 6968.3335
$ws.Range("I65").Value = 6968.3335
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 34841.6675
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31721.6675
$ws.Range("H68").Value = 62500.5
$ws.Range("J68").Value = 62500.5
$ws.Range("L68").Value = 62500.5
$ws.Range("N68").Value = -63998.5
$ws.Range("H71").Value = 62500.5
$ws.Range("J71").Value = 62500.5
$ws.Range("L71").Value = 187501.5
$ws.Range("N71").Value = -194989.5
$ws.Range("H96").Value = 34500
$ws.Range("J96").Value = 34500
$ws.Range("L96").Value = 34500
$ws.Range("N96").Value = -39992
$ws.Range("H99").Value = 2002010
$ws.Range("I99").Value = 2002010
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2002010
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -2000512
$ws.Range("H106").Value = 29000
$ws.Range("J106").Value = 29000
$ws.Range("L106").Value = 29000
$ws.Range("N106").Value = -31524
$ws.Range("H107").Value = 778.625
$ws.Range("I107").Value = 415
$ws.Range("K107").Value = 415
$ws.Range("M107").Value = 1505
$ws.Range("H118").Value = 75000
$ws.Range("J118").Value = 75000
$ws.Range("L118").Value = 75000
$ws.Range("N118").Value = -78314
$ws.Range("H126").Value = 2002010
$ws.Range("I126").Value = 2002010
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6006030
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -6003560
$ws.Range("H132").Value = 4867.5454
$ws.Range("I132").Value = 4867.5454
$ws.Range("K132").Value = 14602.6362
$ws.Range("M132").Value = -12072.6362
$ws.Range("H134").Value = 4380
$ws.Range("I134").Value = 3564.6924
$ws.Range("K134").Value = 10694.0772
$ws.Range("M134").Value = -8159.0772
$ws.Range("H135").Value = 190305.8
$ws.Range("J135").Value = 190305.8
$ws.Range("L135").Value = 190305.8
$ws.Range("N135").Value = -200445.8
$ws.Range("H136").Value = 2266.2173
$ws.Range("I136").Value = 2266.2173
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6798.651899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -4248.651899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 71.64
$ws.Range("I4").Value = 72.07071000000001
$ws.Range("J4").Value = 29
$ws.Range("K4").Value = 216.21213
$ws.Range("L4").Value = 87
$ws.Range("M4").Value = -104.21213
$ws.Range("N4").Value = -311
$ws.Range("H31").Value = 4499
$ws.Range("J31").Value = 4499
$ws.Range("L31").Value = 13497
$ws.Range("N31").Value = -14073
$ws.Range("H114").Value = 832.5
$ws.Range("I114").Value = 400
$ws.Range("J114").Value = 1265
$ws.Range("K114").Value = 1200
$ws.Range("L114").Value = 3795
$ws.Range("M114").Value = 2054
$ws.Range("N114").Value = -10303
$ws.Range("H122").Value = 1983
$ws.Range("I122").Value = 669.5714
$ws.Range("J122").Value = 2557.625
$ws.Range("K122").Value = 6026.1426
$ws.Range("L122").Value = 23018.625
$ws.Range("M122").Value = -3576.1426
$ws.Range("N122").Value = -27918.625
$ws.Range("H133").Value = 4331.5
$ws.Range("I133").Value = 3998
$ws.Range("K133").Value = 11994
$ws.Range("M133").Value = -6934

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9491334
$ws.Range("I11").Value = 15811111
$ws.Range("J11").Value = 11666.667
$ws.Range("K11").Value = 15811111
$ws.Range("L11").Value = 11666.667
$ws.Range("M11").Value = -15810972
$ws.Range("N11").Value = -11944.667
$ws.Range("H12").Value = 933350
$ws.Range("I12").Value = 50
$ws.Range("K12").Value = 50
$ws.Range("M12").Value = 90
$ws.Range("H29").Value = 4000
$ws.Range("I29").Value = 4000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 4000
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -3710
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H42").Value = 12363
$ws.Range("I42").Value = 12363
$ws.Range("K42").Value = 12363
$ws.Range("M42").Value = -11878
$ws.Range("H43").Value = 13402.8
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 15671.333
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 15671.333
$ws.Range("M43").Value = -9849
$ws.Range("N43").Value = -15973.333
$ws.Range("H80").Value = 4962.613
$ws.Range("I80").Value = 4219
$ws.Range("J80").Value = 5659.75
$ws.Range("K80").Value = 4219
$ws.Range("L80").Value = 5659.75
$ws.Range("M80").Value = -3221
$ws.Range("N80").Value = -7655.75
$ws.Range("H83").Value = 4962.613
$ws.Range("I83").Value = 4219
$ws.Range("J83").Value = 5659.75
$ws.Range("K83").Value = 21095
$ws.Range("L83").Value = 28298.75
$ws.Range("M83").Value = -16103
$ws.Range("N83").Value = -38282.75
$ws.Range("H102").Value = 7705.885
$ws.Range("I102").Value = 6516.091
$ws.Range("K102").Value = 6516.091
$ws.Range("M102").Value = -4894.091
$ws.Range("H107").Value = 971.25
$ws.Range("I107").Value = 150
$ws.Range("J107").Value = 1245
$ws.Range("K107").Value = 150
$ws.Range("L107").Value = 1245
$ws.Range("M107").Value = 1770
$ws.Range("N107").Value = -5085
$ws.Range("H113").Value = 25645928
$ws.Range("J113").Value = 9200
$ws.Range("L113").Value = 9200
$ws.Range("N113").Value = -13540
$ws.Range("H115").Value = 12363
$ws.Range("I115").Value = 12363
$ws.Range("K115").Value = 12363
$ws.Range("M115").Value = -11188
$ws.Range("H126").Value = 2800.4375
$ws.Range("I126").Value = 2758.4285
$ws.Range("K126").Value = 8275.2855
$ws.Range("M126").Value = -5805.2855
$ws.Range("H132").Value = 2411.818
$ws.Range("I132").Value = 2215.1428
$ws.Range("J132").Value = 2756
$ws.Range("K132").Value = 6645.428400000001
$ws.Range("L132").Value = 8268
$ws.Range("M132").Value = -4115.428400000001
$ws.Range("N132").Value = -13328

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2934.5454
$ws.Range("I7").Value = 2319.7778
$ws.Range("K7").Value = 2319.7778
$ws.Range("M7").Value = -2207.7778
$ws.Range("H16").Value = 580.4286
$ws.Range("J16").Value = 691.1667
$ws.Range("L16").Value = 691.1667
$ws.Range("N16").Value = -1031.1667
$ws.Range("H19").Value = 4581
$ws.Range("J19").Value = 7633
$ws.Range("L19").Value = 7633
$ws.Range("N19").Value = -7973
$ws.Range("H22").Value = 3096.76
$ws.Range("I22").Value = 2369.25
$ws.Range("K22").Value = 2369.25
$ws.Range("M22").Value = -2074.25
$ws.Range("H27").Value = 3096.76
$ws.Range("I27").Value = 2369.25
$ws.Range("K27").Value = 2369.25
$ws.Range("M27").Value = -2262.25
$ws.Range("H30").Value = 177.33333
$ws.Range("J30").Value = 500
$ws.Range("L30").Value = 500
$ws.Range("N30").Value = -716
$ws.Range("H35").Value = 9511.4
$ws.Range("I35").Value = 519
$ws.Range("J35").Value = 23000
$ws.Range("K35").Value = 519
$ws.Range("L35").Value = 23000
$ws.Range("M35").Value = -183
$ws.Range("N35").Value = -23672
$ws.Range("H40").Value = 4834.387
$ws.Range("I40").Value = 3546
$ws.Range("J40").Value = 7540
$ws.Range("K40").Value = 3546
$ws.Range("L40").Value = 7540
$ws.Range("M40").Value = -3410
$ws.Range("N40").Value = -7812
$ws.Range("H46").Value = 905.6896400000001
$ws.Range("I46").Value = 686.1818
$ws.Range("J46").Value = 1595.5714
$ws.Range("K46").Value = 686.1818
$ws.Range("L46").Value = 1595.5714
$ws.Range("M46").Value = -498.1818
$ws.Range("N46").Value = -1971.5714
$ws.Range("H55").Value = 878.9091
$ws.Range("I55").Value = 147
$ws.Range("J55").Value = 2830.6667
$ws.Range("K55").Value = 147
$ws.Range("L55").Value = 2830.6667
$ws.Range("M55").Value = 26
$ws.Range("N55").Value = -3176.6667
$ws.Range("H68").Value = 2956
$ws.Range("I68").Value = 2956
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2956
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -2207
$ws.Range("H71").Value = 2956
$ws.Range("I71").Value = 2956
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14780
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -11036
$ws.Range("H122").Value = 5128.5
$ws.Range("I122").Value = 5087.533
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 15262.599
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -12812.599
$ws.Range("N122").Value = -20900.0005
$ws.Range("H126").Value = 2934.5454
$ws.Range("I126").Value = 2319.7778
$ws.Range("K126").Value = 6959.3334
$ws.Range("M126").Value = -4489.3334
$ws.Range("H132").Value = 2948.4358
$ws.Range("I132").Value = 3002.1924
$ws.Range("J132").Value = 2840.923
$ws.Range("K132").Value = 9006.5772
$ws.Range("L132").Value = 8522.769
$ws.Range("M132").Value = -6476.5772
$ws.Range("N132").Value = -13582.769

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 7505000
$ws.Range("I29").Value = 15000000
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 15000000
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -14999710
$ws.Range("N29").Value = -10580
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H81").Value = 1477.5
$ws.Range("I81").Value = 1477.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2955
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -1894
$ws.Range("H84").Value = 1477.5
$ws.Range("I84").Value = 1477.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 14775
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -9471
$ws.Range("H113").Value = 934.5862
$ws.Range("I113").Value = 658.7857
$ws.Range("K113").Value = 1976.3571
$ws.Range("M113").Value = 193.6428999999998
$ws.Range("H122").Value = 5888.85
$ws.Range("I122").Value = 3420.4443
$ws.Range("J122").Value = 7908.4546
$ws.Range("K122").Value = 10261.3329
$ws.Range("L122").Value = 23725.3638
$ws.Range("M122").Value = -7811.332900000001
$ws.Range("N122").Value = -28625.3638
$ws.Range("H136").Value = 1790.1538
$ws.Range("I136").Value = 1030.2222
$ws.Range("K136").Value = 3090.6666
$ws.Range("M136").Value = -540.6665999999996
